# Scheduled runner update: refresh market-price / profit figures on the
# Leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 790.95
$ws.Range("I18").Value = 785.2105
$ws.Range("K18").Value = 785.2105
$ws.Range("M18").Value = -501.2105

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1807.8
$ws.Range("I40").Value = 1576.4706
$ws.Range("J40").Value = 2110.3076
$ws.Range("K40").Value = 1576.4706
$ws.Range("L40").Value = 2110.3076
$ws.Range("M40").Value = -1401.4706
$ws.Range("N40").Value = -2460.3076

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2933.5
$ws.Range("I125").Value = 2810.6667
$ws.Range("J125").Value = 3007.2
$ws.Range("K125").Value = 25296.0003
$ws.Range("L125").Value = 27064.8
$ws.Range("M125").Value = -22836.0003
$ws.Range("N125").Value = -31984.8

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2548142.8
$ws.Range("I138").Value = 1515.2142
$ws.Range("J138").Value = 3946291
$ws.Range("K138").Value = 4545.642599999999
$ws.Range("L138").Value = 11838873
$ws.Range("M138").Value = 594.3574000000008
$ws.Range("N138").Value = -11849153

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2063.6072
$ws.Range("I2").Value = 2029.5264
$ws.Range("J2").Value = 2135.5557
$ws.Range("K2").Value = 2029.5264
$ws.Range("L2").Value = 2135.5557
$ws.Range("M2").Value = -1916.5264
$ws.Range("N2").Value = -2361.5557

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2035.1489
$ws.Range("I45").Value = 1666.8125
$ws.Range("J45").Value = 2820.9333
$ws.Range("K45").Value = 1666.8125
$ws.Range("L45").Value = 2820.9333
$ws.Range("M45").Value = -1289.8125
$ws.Range("N45").Value = -3574.9333

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3073.0435
$ws.Range("I63").Value = 2510
$ws.Range("J63").Value = 3373.3333
$ws.Range("K63").Value = 2510
$ws.Range("L63").Value = 3373.3333
$ws.Range("M63").Value = -1824
$ws.Range("N63").Value = -4745.3333

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3073.0435
$ws.Range("I66").Value = 2510
$ws.Range("J66").Value = 3373.3333
$ws.Range("K66").Value = 12550
$ws.Range("L66").Value = 16866.6665
$ws.Range("M66").Value = -9118
$ws.Range("N66").Value = -23730.6665

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2063.6072
$ws.Range("I116").Value = 2029.5264
$ws.Range("J116").Value = 2135.5557
$ws.Range("K116").Value = 2029.5264
$ws.Range("L116").Value = 2135.5557
$ws.Range("M116").Value = 264.4736
$ws.Range("N116").Value = -6723.5557

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2063.6072
$ws.Range("I3").Value = 2029.5264
$ws.Range("J3").Value = 2135.5557
$ws.Range("K3").Value = 2029.5264
$ws.Range("L3").Value = 2135.5557
$ws.Range("M3").Value = -1915.5264
$ws.Range("N3").Value = -2363.5557

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 886.0454999999999
$ws.Range("I80").Value = 1263
$ws.Range("J80").Value = 509.0909
$ws.Range("K80").Value = 1263
$ws.Range("L80").Value = 509.0909
$ws.Range("M80").Value = -265
$ws.Range("N80").Value = -2505.0909

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 886.0454999999999
$ws.Range("I83").Value = 1263
$ws.Range("J83").Value = 509.0909
$ws.Range("K83").Value = 6315
$ws.Range("L83").Value = 2545.4545
$ws.Range("M83").Value = -1323
$ws.Range("N83").Value = -12529.4545

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7577330.5
$ws.Range("I31").Value = 1058.5682
$ws.Range("J31").Value = 22729874
$ws.Range("K31").Value = 1058.5682
$ws.Range("L31").Value = 22729874
$ws.Range("M31").Value = -763.5681999999999
$ws.Range("N31").Value = -22730464

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7577330.5
$ws.Range("I34").Value = 1058.5682
$ws.Range("J34").Value = 22729874
$ws.Range("K34").Value = 1058.5682
$ws.Range("L34").Value = 22729874
$ws.Range("M34").Value = -856.5681999999999
$ws.Range("N34").Value = -22730278

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 3849.5
$ws.Range("I41").Value = 3849.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3849.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -3421.5
$ws.Range("N41").ClearContents()

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 12000
$ws.Range("J60").Value = 12000
$ws.Range("L60").Value = 12000
$ws.Range("N60").Value = -13022

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 825501.5600000001
$ws.Range("I132").Value = 2277.2368
$ws.Range("J132").Value = 5294433.5
$ws.Range("K132").Value = 6831.7104
$ws.Range("L132").Value = 15883300.5
$ws.Range("M132").Value = -4301.7104
$ws.Range("N132").Value = -15888360.5

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4076.6667
$ws.Range("I56").Value = 4076.6667
$ws.Range("K56").Value = 4076.6667
$ws.Range("M56").Value = -3546.6667

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3677.524
$ws.Range("I63").Value = 699.8333
$ws.Range("J63").Value = 4868.6
$ws.Range("K63").Value = 2099.4999
$ws.Range("L63").Value = 14605.8
$ws.Range("M63").Value = -1350.4999
$ws.Range("N63").Value = -16103.8

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5234.5884
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 5499.143
$ws.Range("K64").Value = 12000
$ws.Range("L64").Value = 16497.429
$ws.Range("M64").Value = -11730
$ws.Range("N64").Value = -17037.429

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3677.524
$ws.Range("I66").Value = 699.8333
$ws.Range("J66").Value = 4868.6
$ws.Range("K66").Value = 6298.4997
$ws.Range("L66").Value = 43817.4
$ws.Range("M66").Value = -2554.4997
$ws.Range("N66").Value = -51305.4

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 5234.5884
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 5499.143
$ws.Range("K67").Value = 12000
$ws.Range("L67").Value = 16497.429
$ws.Range("M67").Value = -11064
$ws.Range("N67").Value = -18369.429

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3871.5264
$ws.Range("I87").Value = 2910.6
$ws.Range("J87").Value = 7475
$ws.Range("K87").Value = 8731.799999999999
$ws.Range("L87").Value = 22425
$ws.Range("M87").Value = -7483.799999999999
$ws.Range("N87").Value = -24921

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 3871.5264
$ws.Range("I90").Value = 2910.6
$ws.Range("J90").Value = 7475
$ws.Range("K90").Value = 26195.4
$ws.Range("L90").Value = 67275
$ws.Range("M90").Value = -19955.4
$ws.Range("N90").Value = -79755

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1667131.9
$ws.Range("I107").Value = 298.57144
$ws.Range("J107").Value = 2281228.5
$ws.Range("K107").Value = 895.71432
$ws.Range("L107").Value = 6843685.5
$ws.Range("M107").Value = 1024.28568
$ws.Range("N107").Value = -6847525.5

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1866.6666
$ws.Range("I140").Value = 621.0526
$ws.Range("J140").Value = 6600
$ws.Range("K140").Value = 1863.1578
$ws.Range("L140").Value = 19800
$ws.Range("M140").Value = 3316.8422
$ws.Range("N140").Value = -30160

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4927.6562
$ws.Range("I70").Value = 4655.84
$ws.Range("J70").Value = 5898.4287
$ws.Range("K70").Value = 4655.84
$ws.Range("L70").Value = 5898.4287
$ws.Range("M70").Value = -4385.84
$ws.Range("N70").Value = -6438.4287

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4927.6562
$ws.Range("I73").Value = 4655.84
$ws.Range("J73").Value = 5898.4287
$ws.Range("K73").Value = 4655.84
$ws.Range("L73").Value = 5898.4287
$ws.Range("M73").Value = -3719.84
$ws.Range("N73").Value = -7770.4287

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4355.3877
$ws.Range("I132").Value = 4294.387
$ws.Range("J132").Value = 4460.4443
$ws.Range("K132").Value = 12883.161
$ws.Range("L132").Value = 13381.3329
$ws.Range("M132").Value = -10353.161
$ws.Range("N132").Value = -18441.3329

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 30569.666
$ws.Range("J46").Value = 30569.666
$ws.Range("L46").Value = 30569.666
$ws.Range("N46").Value = -31031.666

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 144656.56
$ws.Range("I122").Value = 23726
$ws.Range("J122").Value = 241401
$ws.Range("K122").Value = 71178
$ws.Range("L122").Value = 724203
$ws.Range("M122").Value = -68728
$ws.Range("N122").Value = -729103

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 30569.666
$ws.Range("J134").Value = 30569.666
$ws.Range("L134").Value = 91708.99800000001
$ws.Range("N134").Value = -96778.99800000001
